$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "Đang sống" (Still alive) in column F ("Năm mất" / Year of death)
# for every data row that does not already have a death year, i.e. rows
# 4-19 and 21-31 (row 20 already has a death year of 2023).
$rows = @(4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,21,22,23,24,25,26,27,28,29,30,31)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "Đang sống"
}

# Update the sheet view: scroll so row 7 is the top-left visible row, and
# select F21:F31 with the active cell at F21.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F21:F31").Select()
